$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $addr, $value) {
    $ws.Range($addr).Value = $value
}

function Clear-Cell($ws, $addr) {
    $ws.Range($addr).ClearContents()
}

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")

# row 17
Set-Cell $ws "H17" 3483.8333
Set-Cell $ws "J17" 3483.8333
Set-Cell $ws "L17" 10451.4999
Set-Cell $ws "N17" -10787.4999

# row 87
Set-Cell $ws "H87" 79995
Set-Cell $ws "J87" 79995
Set-Cell $ws "L87" 79995
Set-Cell $ws "N87" -82491

# row 90
Set-Cell $ws "H90" 79995
Set-Cell $ws "J90" 79995
Set-Cell $ws "L90" 239985
Set-Cell $ws "N90" -252465

# row 106
Set-Cell $ws "H106" 3500
Set-Cell $ws "I106" 3500
Set-Cell $ws "K106" 3500
Set-Cell $ws "M106" -2869

# row 112
Set-Cell $ws "H112" 3172.1155
Set-Cell $ws "J112" 3469.1304
Set-Cell $ws "L112" 10407.3912
Set-Cell $ws "N112" -12623.3912

# row 118 (new M118 cell added)
Set-Cell $ws "H118" 250
Set-Cell $ws "I118" 250
Set-Cell $ws "K118" 750
Set-Cell $ws "M118" 907

# row 135
Set-Cell $ws "H135" 131.42857
Set-Cell $ws "I135" 133.84616
Set-Cell $ws "K135" 1204.61544
Set-Cell $ws "M135" 1330.38456

# row 137
Set-Cell $ws "H137" 1790.6666
Set-Cell $ws "J137" 1500
Set-Cell $ws "L137" 4500
Set-Cell $ws "N137" -9600

# row 141
Set-Cell $ws "H141" 939.56525
Set-Cell $ws "I141" 939.56525
Set-Cell $ws "K141" 2818.69575
Set-Cell $ws "M141" 2361.30425

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")

# row 61
Set-Cell $ws "H61" 2459.4666
Set-Cell $ws "I61" 1489.2
Set-Cell $ws "K61" 1489.2
Set-Cell $ws "M61" -1277.2

# row 74
Set-Cell $ws "H74" 1173.9
Set-Cell $ws "I74" 1180.5555
Set-Cell $ws "K74" 1180.5555
Set-Cell $ws "M74" -306.5554999999999

# row 77
Set-Cell $ws "H77" 1173.9
Set-Cell $ws "I77" 1180.5555
Set-Cell $ws "K77" 5902.7775
Set-Cell $ws "M77" -1534.7775

# row 136
Set-Cell $ws "H136" 2459.4666
Set-Cell $ws "I136" 1489.2
Set-Cell $ws "K136" 4467.6
Set-Cell $ws "M136" -1917.6

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")

# row 86
Set-Cell $ws "H86" 3930.8572
Set-Cell $ws "I86" 2962.6667
Set-Cell $ws "J86" 4657
Set-Cell $ws "K86" 2962.6667
Set-Cell $ws "L86" 4657
Set-Cell $ws "M86" -1839.6667
Set-Cell $ws "N86" -6903

# row 89
Set-Cell $ws "H89" 3930.8572
Set-Cell $ws "I89" 2962.6667
Set-Cell $ws "J89" 4657
Set-Cell $ws "K89" 14813.3335
Set-Cell $ws "L89" 23285
Set-Cell $ws "M89" -9197.333500000001
Set-Cell $ws "N89" -34517

# row 134
Set-Cell $ws "H134" 3751.6667
Set-Cell $ws "I134" 3773.6667
Set-Cell $ws "K134" 11321.0001
Set-Cell $ws "M134" -8786.000100000001

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")

# row 31 (M31 removed)
Set-Cell $ws "H31" 5000
Set-Cell $ws "I31" 0
Set-Cell $ws "K31" 0
Clear-Cell $ws "M31"

# row 34 (M34 removed)
Set-Cell $ws "H34" 5000
Set-Cell $ws "I34" 0
Set-Cell $ws "K34" 0
Clear-Cell $ws "M34"

# row 70 (N70 removed)
Set-Cell $ws "H70" 0
Set-Cell $ws "J70" 0
Set-Cell $ws "L70" 0
Clear-Cell $ws "N70"

# row 73 (N73 removed)
Set-Cell $ws "H73" 0
Set-Cell $ws "J73" 0
Set-Cell $ws "L73" 0
Clear-Cell $ws "N73"

# row 132
Set-Cell $ws "H132" 2756.9
Set-Cell $ws "I132" 2196.375
Set-Cell $ws "K132" 6589.125
Set-Cell $ws "M132" -4059.125

# row 134
Set-Cell $ws "H134" 2844.125
Set-Cell $ws "I134" 2844.125
Set-Cell $ws "K134" 8532.375
Set-Cell $ws "M134" -5997.375

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")

# row 2
Set-Cell $ws "H2" 102.30882
Set-Cell $ws "I2" 75.5
Set-Cell $ws "J2" 103.121216
Set-Cell $ws "K2" 453
Set-Cell $ws "L2" 618.727296
Set-Cell $ws "M2" -340
Set-Cell $ws "N2" -844.727296

# row 8
Set-Cell $ws "H8" 932.125
Set-Cell $ws "I8" 932.125
Set-Cell $ws "K8" 2796.375
Set-Cell $ws "M8" -2657.375

# row 75
Set-Cell $ws "H75" 500
Set-Cell $ws "I75" 500
Set-Cell $ws "K75" 1500
Set-Cell $ws "M75" -502

# row 78
Set-Cell $ws "H78" 500
Set-Cell $ws "I78" 500
Set-Cell $ws "K78" 4500
Set-Cell $ws "M78" 492

# row 107 (M107 removed)
Set-Cell $ws "H107" 65
Set-Cell $ws "I107" 0
Set-Cell $ws "J107" 65
Set-Cell $ws "K107" 0
Set-Cell $ws "L107" 195
Clear-Cell $ws "M107"
Set-Cell $ws "N107" -4035

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")

# row 105
Set-Cell $ws "H105" 670999.5
Set-Cell $ws "J105" 670999.5
Set-Cell $ws "L105" 670999.5
Set-Cell $ws "N105" -677987.5

# row 113 (N113 removed)
Set-Cell $ws "H113" 3733
Set-Cell $ws "I113" 3733
Set-Cell $ws "J113" 0
Set-Cell $ws "K113" 3733
Set-Cell $ws "L113" 0
Set-Cell $ws "M113" -1563
Clear-Cell $ws "N113"

# row 132
Set-Cell $ws "H132" 2531.875
Set-Cell $ws "I132" 1959.4166
Set-Cell $ws "J132" 4249.25
Set-Cell $ws "K132" 5878.2498
Set-Cell $ws "L132" 12747.75
Set-Cell $ws "M132" -3348.2498
Set-Cell $ws "N132" -17807.75

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")

# row 136 (new M136 cell added)
Set-Cell $ws "H136" 8999.5
Set-Cell $ws "I136" 8999
Set-Cell $ws "J136" 9000
Set-Cell $ws "K136" 26997
Set-Cell $ws "L136" 27000
Set-Cell $ws "M136" -24447
Set-Cell $ws "N136" -32100

# ---------------- WVR ----------------
$ws = $wb.Worksheets.Item("WVR")

# row 113
Set-Cell $ws "H113" 893.2308
Set-Cell $ws "I113" 717.4167
Set-Cell $ws "J113" 3003
Set-Cell $ws "K113" 2152.2501
Set-Cell $ws "L113" 9009
Set-Cell $ws "M113" 17.7498999999998
Set-Cell $ws "N113" -13349

# row 124 (new N124 cell added)
Set-Cell $ws "H124" 100000
Set-Cell $ws "J124" 100000
Set-Cell $ws "L124" 100000
Set-Cell $ws "N124" -109820

# row 136
Set-Cell $ws "H136" 1549.1428
Set-Cell $ws "I136" 1549.1428
Set-Cell $ws "K136" 4647.428400000001
Set-Cell $ws "M136" -2097.428400000001
